$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Summary header updates (totals recomputed after adding new worker rows)
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1536622      # VALOR MORA (total)
$ws.Range("C13").Value = 17           # Cant. Trabajadores
$ws.Range("F13").Value = 11           # Cant. Periodos

# ---------------------------------------------------------------------------
# 2. Relocate the signature-block footer from rows 34-35 down to rows 48-49
#    (keeps the same gap below the now-longer detail table). This must run
#    before the detail rows 30-43 are populated, since the footer currently
#    overlaps that range.
# ---------------------------------------------------------------------------
$ws.Range("B34:C34").UnMerge() | Out-Null
$ws.Range("B35:C35").UnMerge() | Out-Null
$ws.Range("H34:J34").UnMerge() | Out-Null
$ws.Range("H35:J35").UnMerge() | Out-Null

$ws.Range("B34:J35").Cut($ws.Range("B48")) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B48:C48").Merge() | Out-Null
$ws.Range("B49:C49").Merge() | Out-Null
$ws.Range("H48:J48").Merge() | Out-Null
$ws.Range("H49:J49").Merge() | Out-Null

# ---------------------------------------------------------------------------
# 3. Prepare the detail table area (rows 16-43) so formatting matches the
#    rest of the table before we pour in the refreshed data.
#    Row 29 used to be the last (bottom-bordered) row of the table; it is no
#    longer last, so give it ordinary formatting first, then stamp the
#    "normal" row format onto all the newly-needed rows 30-43, and finally
#    re-apply the special bottom-bordered "last row" formatting to the new
#    true last row (43).
# ---------------------------------------------------------------------------
$ws.Range("B28:J28").Copy() | Out-Null
$ws.Range("B29:J29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B28:J28").Copy() | Out-Null
$ws.Range("B30:J43").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B29:J29").Copy() | Out-Null
$ws.Range("B43:J43").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (bottom border)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Refresh the detail table content (rows 16-43)
# ---------------------------------------------------------------------------
$rows = @(
    @{Row=16; Doc='1050949292'; Name='LINA MARIA CASTELLON ESPINOSA'; Period='1709'; F=30679; G=766979},
    @{Row=17; Doc='1050949292'; Name='LINA MARIA CASTELLON ESPINOSA'; Period='1710'; F=30679; G=766979},
    @{Row=18; Doc='1002241486'; Name='KATHERINE DEL CARMEN PAJARO MENDOZA'; Period='2402'; F=7936; G=1487922},
    @{Row=19; Doc='1002241486'; Name='KATHERINE DEL CARMEN PAJARO MENDOZA'; Period='2403'; F=59516; G=1487922},
    @{Row=20; Doc='1002241486'; Name='KATHERINE DEL CARMEN PAJARO MENDOZA'; Period='2404'; F=59516; G=1487922},
    @{Row=21; Doc='1042579787'; Name='EVERLIN VANESA AGAMEZ MARRUGO'; Period='2411'; F=45500; G=1423500},
    @{Row=22; Doc='1042579787'; Name='EVERLIN VANESA AGAMEZ MARRUGO'; Period='2412'; F=54600; G=1423500},
    @{Row=23; Doc='1051417343'; Name='WILSON DANOVIS MUÑOZ MIRANDA'; Period='2501'; F=1721; G=1850000},
    @{Row=24; Doc='1042579787'; Name='EVERLIN VANESA AGAMEZ MARRUGO'; Period='2501'; F=54600; G=1423500},
    @{Row=25; Doc='1051417343'; Name='WILSON DANOVIS MUÑOZ MIRANDA'; Period='2502'; F=7375; G=1850000},
    @{Row=26; Doc='1042579787'; Name='EVERLIN VANESA AGAMEZ MARRUGO'; Period='2502'; F=54600; G=1423500},
    @{Row=27; Doc='45748640'; Name='MARIA DEL ROSARIO PARRA TORRES'; Period='2507'; F=37000; G=1850000},
    @{Row=28; Doc='3838464'; Name='IVAN DAVID MACARENO VERGARA'; Period='2507'; F=66600; G=1850000},
    @{Row=29; Doc='8853279'; Name='SIR JAVIER HERNANDEZ JIMENEZ'; Period='2508'; F=74000; G=1850000},
    @{Row=30; Doc='45748640'; Name='MARIA DEL ROSARIO PARRA TORRES'; Period='2508'; F=74000; G=1850000},
    @{Row=31; Doc='8950085'; Name='NAFER ENRIQUE CARO CERPA'; Period='2508'; F=74000; G=1850000},
    @{Row=32; Doc='9098296'; Name='ELKIN MALLARINO LLERENA'; Period='2508'; F=74000; G=1850000},
    @{Row=33; Doc='1052068040'; Name='LINA MARGARITA MARTINEZ ROBLES'; Period='2508'; F=75600; G=1890000},
    @{Row=34; Doc='19897222'; Name='ORLANDO ALFARO PARRA'; Period='2508'; F=74000; G=1850000},
    @{Row=35; Doc='1051417343'; Name='WILSON DANOVIS MUÑOZ MIRANDA'; Period='2508'; F=74000; G=1850000},
    @{Row=36; Doc='1051417083'; Name='GUSTAVO ADOLFO JULIO REBOLLEDO'; Period='2508'; F=74000; G=1850000},
    @{Row=37; Doc='1143333700'; Name='STEPHANI DEL CARMEN VALENCIA OROZCO'; Period='2508'; F=56940; G=1423500},
    @{Row=38; Doc='1051419682'; Name='ANDRES FELIPE ARRIETA DAZA'; Period='2508'; F=74000; G=1850000},
    @{Row=39; Doc='1192763716'; Name='MARGELIS ARROYO ZAMBRANO'; Period='2508'; F=56940; G=1423500},
    @{Row=40; Doc='1047485369'; Name='HARTLEY PRETTEL GALAN'; Period='2508'; F=56940; G=1423500},
    @{Row=41; Doc='1002196824'; Name='MELISSA PAOLA MORENO DAVILA'; Period='2508'; F=56940; G=1423500},
    @{Row=42; Doc='3838464'; Name='IVAN DAVID MACARENO VERGARA'; Period='2508'; F=74000; G=1850000},
    @{Row=43; Doc='1042579787'; Name='EVERLIN VANESA AGAMEZ MARRUGO'; Period='2508'; F=56940; G=1423500}
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 2).Value = "CC"
    $ws.Cells.Item($n, 3).Value = $r.Doc
    $ws.Cells.Item($n, 4).Value = $r.Name
    $ws.Cells.Item($n, 5).Value = $r.Period
    $ws.Cells.Item($n, 6).Value = $r.F
    $ws.Cells.Item($n, 7).Value = $r.G
}
